$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing C and D columns (rows 2..28) before overwriting them,
# then write them back in reverse row order. Also convert column B from the
# shared string "1000" to the numeric value 1000.

$firstRow = 2
$lastRow = 28

$cVals = @()
$dVals = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cVals += $ws.Cells.Item($r, 3).Value2
    $dVals += $ws.Cells.Item($r, 4).Value2
}

$n = $cVals.Count
for ($i = 0; $i -lt $n; $i++) {
    $r = $firstRow + $i
    $srcIndex = $n - 1 - $i
    $ws.Cells.Item($r, 2).Value2 = 1000
    $ws.Cells.Item($r, 3).Value2 = $cVals[$srcIndex]
    $ws.Cells.Item($r, 4).Value2 = $dVals[$srcIndex]
}

# Update the selection to mirror the author's last active cell.
$ws.Range("J9").Select()
